$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update N[0] count
$ws.Range("C36").Value = 1

# Update N[1] count
$ws.Range("C37").Value = 21

# Update increment (deltaN)
$ws.Range("C38").Value = 5

# Add comment about period in G37, entered with leading apostrophe so Excel
# treats it as literal text (quotePrefix) rather than a formula
$ws.Range("G37").Value = "'= 12 sek / U"

# Update the active selection to C38
$ws.Range("C38").Select()
